$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''88.893.67'
$ws.Range("E2").Value = '  -0.06%  '
$ws.Range("D3").Value = '''3.141.55'
$ws.Range("E3").Value = '  -3.97%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''208.33'
$ws.Range("E5").Value = '  -1.43%  '
$ws.Range("D6").Value = '''609.28'
$ws.Range("E6").Value = '  -2.52%  '
$ws.Range("D7").Value = '''0.378'
$ws.Range("E7").Value = '  +1.43%  '
$ws.Range("D8").Value = '''0.677'
$ws.Range("E8").Value = '  -4.52%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").Value = '''3.138.99'
$ws.Range("E10").Value = '  -3.92%  '
$ws.Range("D11").Value = '''0.564'
$ws.Range("E11").Value = '  -1.27%  '
$ws.Range("D12").Value = '''0.175'
$ws.Range("E12").Value = '  -6.35%  '
$ws.Range("D13").Value = '''0.0000247'
$ws.Range("E13").Value = '  -5.03%  '
$ws.Range("D14").Value = '''88.986.94'
$ws.Range("E14").Value = '  -0.02%  '
$ws.Range("D15").Value = '''3.710.34'
$ws.Range("E15").Value = '  -4.30%  '
$ws.Range("D16").Value = '''5.15'
$ws.Range("E16").Value = '  -5.13%  '
$ws.Range("D17").Value = '''32.04'
$ws.Range("E17").Value = '  -5.38%  '
$ws.Range("D18").Value = '''3.172.02'
$ws.Range("E18").Value = '  -3.99%  '
$ws.Range("D19").Value = '''3.20'
$ws.Range("E19").Value = '  +3.42%  '
$ws.Range("D20").Value = '''13.19'
$ws.Range("E20").Value = '  -5.75%  '
$ws.Range("D21").Value = '''430.75'
$ws.Range("E21").Value = '  -1.02%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '''8.42'
$ws.Range("E22").Value = '  -4.56%  '
$ws.Range("B23").Value = 'PEPE'
$ws.Range("C23").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D23").Value = '''0.0000182'
$ws.Range("E23").Value = '  +34.86%  '
$ws.Range("D24").Value = '''5.01'
$ws.Range("E24").Value = '  -5.70%  '
$ws.Range("D25").Value = '''5.03'
$ws.Range("E25").Value = '  -3.43%  '
$ws.Range("D26").Value = '''11.48'
$ws.Range("E26").Value = '  -5.42%  '
$ws.Range("D27").Value = '''3.353.69'
$ws.Range("E27").Value = '  -3.13%  '
$ws.Range("D28").Value = '''73.80'
$ws.Range("E28").Value = '  -3.71%  '
$ws.Range("D30").Value = '''0.164'
$ws.Range("E30").Value = '  -8.25%  '
$ws.Range("D31").Value = '''1.00'
$ws.Range("E31").Value = '  +0.33%  '
$ws.Range("D32").Value = '''3.96'
$ws.Range("E32").Value = '  +28.61%  '
$ws.Range("D33").Value = '''8.27'
$ws.Range("E33").Value = '  -5.97%  '
$ws.Range("D34").Value = '''522.03'
$ws.Range("E34").Value = '  -6.69%  '
$ws.Range("D35").Value = '''6.81'
$ws.Range("E35").Value = '  -3.44%  '
$ws.Range("D36").Value = '''1.84'
$ws.Range("E36").Value = '  -5.85%  '
$ws.Range("D37").Value = '''1.24'
$ws.Range("E37").Value = '  -8.69%  '
$ws.Range("D38").Value = '''22.25'
$ws.Range("E38").Value = '  +1.88%  '
$ws.Range("D39").Value = '''21.56'
$ws.Range("E39").Value = '  -4.71%  '
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").Value = '''0.999'
$ws.Range("E40").Value = '  -0.08%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = '''0.126'
$ws.Range("E41").Value = '  -9.94%  '
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("D43").Value = '''0.368'
$ws.Range("E43").Value = '  -8.02%  '
$ws.Range("D44").Value = '''1.87'
$ws.Range("E44").Value = '  -7.41%  '
$ws.Range("D45").Value = '''149.42'
$ws.Range("E45").Value = '  -3.96%  '
$ws.Range("D46").Value = '''44.06'
$ws.Range("E46").Value = '  -1.74%  '
$ws.Range("D47").Value = '''169.19'
$ws.Range("E47").Value = '  -5.96%  '
$ws.Range("D48").Value = '''0.122'
$ws.Range("E48").Value = '  -9.01%  '
$ws.Range("D49").Value = '''1.21'
$ws.Range("E49").Value = '  -7.07%  '
$ws.Range("B50").Value = 'ARBITRUM'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D50").Value = '''0.598'
$ws.Range("E50").Value = '  -3.98%  '
$ws.Range("B51").Value = 'Filecoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D51").Value = '''3.98'
$ws.Range("E51").Value = '  -5.24%  '
